$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "date" and "time" field types used by the JGI follow_arrival form are
# being changed to "text" (per the commit message: "Changed date and time
# fields for JGI app"). This affects three rows in the survey sheet:
#   - row 2  (FA_FOL_date):   date -> text
#   - row 9  (FA_time_start): time -> text
#   - row 10 (FA_time_end):   time -> text
$ws.Range("C2").Value = "text"
$ws.Range("C9").Value = "text"
$ws.Range("C10").Value = "text"

# Move/restore the active selection to C11 (matches the saved cursor
# position after the edit).
$ws.Range("C11").Select()
